$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab
$ws.Name = "Super Handball League"

# Update team names (column B) per row - lookup fix
$ws.Cells.Item(2, 2).Value = "KRAS/Volendam HS1"
$ws.Cells.Item(3, 2).Value = "Green Park/Handbal Aalsmeer HS1"
$ws.Cells.Item(4, 2).Value = "Herpertz/Bevo HC HS1"
$ws.Cells.Item(5, 2).Value = "HC Visé BM HS1"
$ws.Cells.Item(6, 2).Value = "KTSV Eupen HS1"
$ws.Cells.Item(7, 2).Value = "Sporting Pelt HS1"
$ws.Cells.Item(8, 2).Value = "HUBO Handbal HS1"
$ws.Cells.Item(9, 2).Value = "Sezoens Achilles Bocholt HS1"
$ws.Cells.Item(10, 2).Value = "JD Techniek/ Hurry-up HS1"
$ws.Cells.Item(11, 2).Value = "LIMBURG LIONS/ Sittardia HS1"
$ws.Cells.Item(12, 2).Value = "Biobest/ Sasja HC HS1"
$ws.Cells.Item(13, 2).Value = "LvanRaak Milieu/Handbal Houten HS1"

# Update standings statistics (columns C-K)
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 41
$ws.Cells.Item(2, 9).Value = 23
$ws.Cells.Item(2, 10).Value = 18
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 36
$ws.Cells.Item(3, 9).Value = 28
$ws.Cells.Item(3, 10).Value = 8
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 31
$ws.Cells.Item(4, 9).Value = 26
$ws.Cells.Item(4, 10).Value = 5
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 38
$ws.Cells.Item(5, 9).Value = 34
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 2
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 28
$ws.Cells.Item(6, 9).Value = 27
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 23
$ws.Cells.Item(7, 9).Value = 22
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 27
$ws.Cells.Item(8, 9).Value = 28
$ws.Cells.Item(8, 10).Value = -1
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 22
$ws.Cells.Item(9, 9).Value = 23
$ws.Cells.Item(9, 10).Value = -1
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 34
$ws.Cells.Item(10, 9).Value = 38
$ws.Cells.Item(10, 10).Value = -4
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 26
$ws.Cells.Item(11, 9).Value = 31
$ws.Cells.Item(11, 10).Value = -5
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 28
$ws.Cells.Item(12, 9).Value = 36
$ws.Cells.Item(12, 10).Value = -8
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 23
$ws.Cells.Item(13, 9).Value = 41
$ws.Cells.Item(13, 10).Value = -18
$ws.Cells.Item(13, 11).Value = 0
